$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 52, shifting existing rows 52..113 down to 53..114.
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new record.
$ws.Cells.Item(52, 1).Value = 7
$ws.Cells.Item(52, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(52, 3).Value = "Ñuble"
$ws.Cells.Item(52, 4).Value = 44902
$ws.Cells.Item(52, 5).Value = 16
$ws.Cells.Item(52, 6).Value = 100112021
$ws.Cells.Item(52, 7).Value = "Ají"
$ws.Cells.Item(52, 8).Value = "Americana (o)"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 60
$ws.Cells.Item(52, 11).Value = 16000
$ws.Cells.Item(52, 12).Value = 17000
$ws.Cells.Item(52, 13).Value = 16500
$ws.Cells.Item(52, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(52, 15).Value = "Región del Maule"
$ws.Cells.Item(52, 16).Value = 1100
$ws.Cells.Item(52, 17).Value = 15
$ws.Cells.Item(52, 18).Value = "Hortaliza"

# Match the date-style formatting used by the rest of column D.
$ws.Cells.Item(52, 4).NumberFormat = $ws.Cells.Item(53, 4).NumberFormat
